$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlTextString = "@"

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Use an existing per-quarter sheet as the style donor so the new sheet's
# header row / index column pick up the same look (bold + border + centered).
$styleDonor = $wb.Worksheets.Item(5)

$styleDonor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

$styleDonor.Range("A2").Copy()
$q1.Range("A2:A8").PasteSpecial($xlPasteFormats)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Columns B, C, D, E, F, G are stored as text (even the numeric-looking
# ones), so force text format before assigning them.
$q1.Range("B2:G8").NumberFormat = $xlTextString

$q1Data = @(
    @(0, "008283", "易方达金融行业股票",                         "34.73", "91.49", "3.53", "1.2260", 9),
    @(1, "501025", "鹏华港股通中证香港银行投资指数（LOF）A",     "9.81",  "94.47", "5.20", "0.5101", 7),
    @(2, "010365", "鹏华港股通中证香港银行投资指数（LOF）C",     "6.07",  "94.47", "5.20", "0.3156", 7),
    @(3, "006809", "泰康港股通中证香港银行投资指数A",             "1.99",  "94.73", "5.19", "0.1033", 7),
    @(4, "010204", "中银港股通优势成长股票",                     "3.19",  "83.00", "3.08", "0.0983", 9),
    @(5, "006810", "泰康港股通中证香港银行投资指数C",             "0.90",  "94.73", "5.19", "0.0467", 7),
    @(6, "162416", "华宝港股通恒生香港35指数(LOF)",              "0.21",  "94.50", "3.29", "0.0069", 10)
)

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $row = $i + 2
    $rec = $q1Data[$i]
    $q1.Cells.Item($row, 1).Value = $rec[0]
    $q1.Cells.Item($row, 2).Value = $rec[1]
    $q1.Cells.Item($row, 3).Value = $rec[2]
    $q1.Cells.Item($row, 4).Value = $rec[3]
    $q1.Cells.Item($row, 5).Value = $rec[4]
    $q1.Cells.Item($row, 6).Value = $rec[5]
    $q1.Cells.Item($row, 7).Value = $rec[6]
    $q1.Cells.Item($row, 8).Value = $rec[7]
}

# The explicit text NumberFormat left a dedicated style behind; the source
# data (like the other per-quarter sheets) keeps these value cells on the
# default (unstyled) format, so strip it back off now that the text type
# has been locked in.
$q1.Range("B2:G8").ClearFormats()

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 7
$total.Cells.Item(2,4).Value = 2.31

# Renumber the index column (A) for every data row, 0..5 top to bottom.
for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/selection (adding sheets shifts
#    the active tab onto the newly inserted one).
# ---------------------------------------------------------------------
$origActive = $wb.Worksheets.Item(1)
$origActive.Activate() | Out-Null
$origActive.Range("A1").Select() | Out-Null

